$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range('E2').Value = '2026-02-07 16:17:56'
$ws.Range('K2').Value = '8.5 MJ/m2'
$ws.Range('E3').Value = '2026-02-07 16:17:58'
$ws.Range('K3').Value = '14.0 MJ/m2'
$ws.Range('L3').Value = '31.7 km/h - 115º 15:49 TU'
$ws.Range('E4').Value = '2026-02-07 16:18:01'
$ws.Range('H4').Value = '''54%'
$ws.Range('C4').Copy() | Out-Null
$ws.Range('H4').PasteSpecial(-4122) | Out-Null
$ws.Range('K4').Value = '10.2 MJ/m2'
$ws.Range('O4').Value = '12.1 °C'
$ws.Range('E5').Value = '2026-02-07 16:18:04'
$ws.Range('K5').Value = '11.8 MJ/m2'
$ws.Range('O5').Value = '10.5 °C'
$ws.Range('E6').Value = '2026-02-07 16:18:07'
$ws.Range('K6').Value = '12.1 MJ/m2'
$ws.Range('O6').Value = '13.1 °C'
$ws.Range('E7').Value = '2026-02-07 16:18:10'
$ws.Range('H7').Value = '''62%'
$ws.Range('C7').Copy() | Out-Null
$ws.Range('H7').PasteSpecial(-4122) | Out-Null
$ws.Range('K7').Value = '12.6 MJ/m2'
$ws.Range('O7').Value = '9.3 °C'
$ws.Range('E8').Value = '2026-02-07 16:18:12'
$ws.Range('H8').Value = '''75%'
$ws.Range('C8').Copy() | Out-Null
$ws.Range('H8').PasteSpecial(-4122) | Out-Null
$ws.Range('K8').Value = '11.5 MJ/m2'
$ws.Range('O8').Value = '8.8 °C'
$ws.Range('E9').Value = '2026-02-07 16:18:15'
$ws.Range('H9').Value = '''86%'
$ws.Range('C9').Copy() | Out-Null
$ws.Range('H9').PasteSpecial(-4122) | Out-Null
$ws.Range('O9').Value = '3.7 °C'
$ws.Range('E10').Value = '2026-02-07 16:18:17'
$ws.Range('O10').Value = '10.4 °C'
$ws.Range('E11').Value = '2026-02-07 16:18:20'
$ws.Range('H11').Value = '''84%'
$ws.Range('C11').Copy() | Out-Null
$ws.Range('H11').PasteSpecial(-4122) | Out-Null
$ws.Range('K11').Value = '8.6 MJ/m2'
$ws.Range('O11').Value = '3.3 °C'
$ws.Range('E12').Value = '2026-02-07 16:18:23'
$ws.Range('K12').Value = '12.3 MJ/m2'
$ws.Range('E13').Value = '2026-02-07 16:18:26'
$ws.Range('O13').Value = '11.4 °C'
$ws.Range('E14').Value = '2026-02-07 16:18:28'
$ws.Range('H14').Value = '''61%'
$ws.Range('C14').Copy() | Out-Null
$ws.Range('H14').PasteSpecial(-4122) | Out-Null
$ws.Range('K14').Value = '8.4 MJ/m2'
$ws.Range('E15').Value = '2026-02-07 16:18:31'
$ws.Range('H15').Value = '''71%'
$ws.Range('C15').Copy() | Out-Null
$ws.Range('H15').PasteSpecial(-4122) | Out-Null
$ws.Range('K15').Value = '11.4 MJ/m2'
$ws.Range('O15').Value = '9.7 °C'
$ws.Range('E16').Value = '2026-02-07 16:18:34'
$ws.Range('K16').Value = '6.6 MJ/m2'
$ws.Range('O16').Value = '3.9 °C'
$ws.Range('E17').Value = '2026-02-07 16:18:36'
$ws.Range('H17').Value = '''86%'
$ws.Range('C17').Copy() | Out-Null
$ws.Range('H17').PasteSpecial(-4122) | Out-Null
$ws.Range('K17').Value = '8.7 MJ/m2'
$ws.Range('O17').Value = '4.9 °C'
$ws.Range('E18').Value = '2026-02-07 16:18:39'
$ws.Range('K18').Value = '5.9 MJ/m2'
$ws.Range('O18').Value = '-5.6 °C'
$ws.Range('E19').Value = '2026-02-07 16:18:42'
$ws.Range('H19').Value = '''81%'
$ws.Range('C19').Copy() | Out-Null
$ws.Range('H19').PasteSpecial(-4122) | Out-Null
$ws.Range('K19').Value = '11.7 MJ/m2'
$ws.Range('O19').Value = '7.0 °C'
$ws.Range('E20').Value = '2026-02-07 16:18:44'
$ws.Range('K20').Value = '10.3 MJ/m2'
$ws.Range('E21').Value = '2026-02-07 16:18:46'
$ws.Range('H21').Value = '''67%'
$ws.Range('C21').Copy() | Out-Null
$ws.Range('H21').PasteSpecial(-4122) | Out-Null
$ws.Range('K21').Value = '11.4 MJ/m2'
$ws.Range('O21').Value = '8.4 °C'
$ws.Range('E22').Value = '2026-02-07 16:18:49'
$ws.Range('K22').Value = '12.7 MJ/m2'
$ws.Range('O22').Value = '10.4 °C'
$ws.Range('E23').Value = '2026-02-07 16:18:52'
$ws.Range('H23').Value = '''76%'
$ws.Range('C23').Copy() | Out-Null
$ws.Range('H23').PasteSpecial(-4122) | Out-Null
$ws.Range('J23').Value = '1003.4 hPa'
$ws.Range('K23').Value = '10.0 MJ/m2'
$ws.Range('O23').Value = '10.5 °C'
$ws.Range('E24').Value = '2026-02-07 16:18:54'
$ws.Range('J24').Value = '1002.9 hPa'
$ws.Range('K24').Value = '9.8 MJ/m2'
$ws.Range('O24').Value = '11.2 °C'
$ws.Range('E25').Value = '2026-02-07 16:18:57'
$ws.Range('K25').Value = '7.1 MJ/m2'
$ws.Range('O25').Value = '2.2 °C'
$ws.Range('E26').Value = '2026-02-07 16:19:00'
$ws.Range('K26').Value = '10.5 MJ/m2'
$ws.Range('E27').Value = '2026-02-07 16:19:03'
$ws.Range('H27').Value = '''75%'
$ws.Range('C27').Copy() | Out-Null
$ws.Range('H27').PasteSpecial(-4122) | Out-Null
$ws.Range('J27').Value = '1003.5 hPa'
$ws.Range('K27').Value = '10.9 MJ/m2'
$ws.Range('O27').Value = '11.5 °C'
$ws.Range('E28').Value = '2026-02-07 16:19:05'
$ws.Range('H28').Value = '''81%'
$ws.Range('C28').Copy() | Out-Null
$ws.Range('H28').PasteSpecial(-4122) | Out-Null
$ws.Range('J28').Value = '1005.5 hPa'
$ws.Range('O28').Value = '4.6 °C'
$ws.Range('E29').Value = '2026-02-07 16:19:08'
$ws.Range('K29').Value = '12.3 MJ/m2'
$ws.Range('E30').Value = '2026-02-07 16:19:11'
$ws.Range('E31').Value = '2026-02-07 16:19:13'
$ws.Range('H31').Value = '''86%'
$ws.Range('C31').Copy() | Out-Null
$ws.Range('H31').PasteSpecial(-4122) | Out-Null
$ws.Range('E32').Value = '2026-02-07 16:19:16'
$ws.Range('K32').Value = '11.8 MJ/m2'
$ws.Range('E33').Value = '2026-02-07 16:19:19'
$ws.Range('H33').Value = '''76%'
$ws.Range('C33').Copy() | Out-Null
$ws.Range('H33').PasteSpecial(-4122) | Out-Null
$ws.Range('O33').Value = '10.4 °C'
$ws.Range('E34').Value = '2026-02-07 16:19:21'
$ws.Range('H34').Value = '''71%'
$ws.Range('C34').Copy() | Out-Null
$ws.Range('H34').PasteSpecial(-4122) | Out-Null
$ws.Range('K34').Value = '10.0 MJ/m2'
$ws.Range('O34').Value = '7.7 °C'
$ws.Range('E35').Value = '2026-02-07 16:19:24'
$ws.Range('K35').Value = '5.8 MJ/m2'
$ws.Range('O35').Value = '-4.3 °C'
$ws.Range('E36').Value = '2026-02-07 16:19:27'
$ws.Range('H36').Value = '''78%'
$ws.Range('C36').Copy() | Out-Null
$ws.Range('H36').PasteSpecial(-4122) | Out-Null
$ws.Range('K36').Value = '11.4 MJ/m2'
$ws.Range('O36').Value = '8.5 °C'
$excel.CutCopyMode = $false
